$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21
$ws.Cells.Item($row, 1).Value = 42625.883379629631
$ws.Cells.Item($row, 2).Value = 22
$ws.Cells.Item($row, 3).Value = 57
$ws.Cells.Item($row, 4).Value = 40
$ws.Cells.Item($row, 5).Value = 100
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 18760
$ws.Cells.Item($row, 8).Value = 16061
$ws.Cells.Item($row, 9).Value = 893
$ws.Cells.Item($row, 10).Value = 214
$ws.Cells.Item($row, 11).Value = 151
$ws.Cells.Item($row, 12).Value = 20
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "Named"
